# Preliminary check-in: rename the ODK "form_id" setting to "table_id" so that
# definitions.csv / properties.csv can be generated from it, and add a new
# "properties" sheet that will be processed into properties.csv.

$wb = $excel.ActiveWorkbook

# --- settings sheet: rename the "form_id" setting row to "table_id" -------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "table_id"

# Move the settings sheet's selection off of the old cursor position.
$settings.Range("A3").Select() | Out-Null

# --- survey sheet: a couple of row-height tweaks that came along with the
#     resave (content is untouched) -----------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Rows.Item(9).RowHeight = 31.45

# --- add the new "properties" sheet after "settings" -----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$props = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$props.Name = "properties"

# Header row
$props.Range("A1").Value = "partition"
$props.Range("B1").Value = "aspect"
$props.Range("C1").Value = "key"
$props.Range("D1").Value = "type"
$props.Range("E1").Value = "value"

# Data row describing the table's column order
$props.Range("A2").Value = "Table"
$props.Range("B2").Value = "default"
$props.Range("C2").Value = "colOrder"
$props.Range("D2").Value = "array"
$props.Range("E2").Value = '["FB_FOL_date","FB_FOL_B_AnimID","FB_begin_feed_time","FB_end_feed_time","FB_duration","FB_FPL_local_food_part","FB_FL_local_food_name","FB_FPL_local_food_part2","FB_FL_local_food_name2","FB_local_food_part_written","FB_local_food_name_written"]'

# Make "properties" the active/selected sheet+cell, matching the new workbook
# state (activeTab points at the properties sheet).
$props.Select() | Out-Null
$props.Range("E4").Select() | Out-Null
